# Insert a new data row at row 249 (shifting existing rows 249-305 down to 250-306)
# and populate it with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 249; all rows from 249 downward shift by one.
$ws.Rows.Item(249).Insert()

# Populate the newly inserted row 249 with the new data record.
$ws.Cells.Item(249, 1).Value  = 3
$ws.Cells.Item(249, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(249, 3).Value  = "Coquimbo"
$ws.Cells.Item(249, 4).Value  = 44943
$ws.Cells.Item(249, 5).Value  = 5
$ws.Cells.Item(249, 6).Value  = "Fruta"
$ws.Cells.Item(249, 7).Value  = 100101
$ws.Cells.Item(249, 8).Value  = "Berries"
$ws.Cells.Item(249, 9).Value  = 100101001
$ws.Cells.Item(249, 10).Value = "Arándano (blue)"
$ws.Cells.Item(249, 11).Value = "Sin especificar"
$ws.Cells.Item(249, 12).Value = "Primera"
$ws.Cells.Item(249, 13).Value = 35
$ws.Cells.Item(249, 14).Value = 4000
$ws.Cells.Item(249, 15).Value = 4000
$ws.Cells.Item(249, 16).Value = 4000
$ws.Cells.Item(249, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(249, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(249, 19).Value = 2000
$ws.Cells.Item(249, 20).Value = 2

# Apply the date number-format (same format used by other "Fecha" cells) to the new D249 cell.
$ws.Cells.Item(249, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
